$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber RowID values in rows 25-29 (62..66 -> 68..72)
$ws.Range("A25").Value = 68
$ws.Range("A26").Value = 69
$ws.Range("A27").Value = 70
$ws.Range("A28").Value = 71
$ws.Range("A29").Value = 72

# Clear the stray IsForeignKey/ReferencedTable/ReferencedColumn data that had
# leaked onto the Inspections rows (InspectorId / MaintenanceActionId) - rows 58-59
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = ""
$ws.Range("M58").Value = ""

$ws.Range("K59").Value = 0
$ws.Range("L59").Value = ""
$ws.Range("M59").Value = ""

# Renumber RowID values in rows 68-73 (67..72 -> 62..67)
$ws.Range("A68").Value = 62
$ws.Range("A69").Value = 63
$ws.Range("A70").Value = 64
$ws.Range("A71").Value = 65
$ws.Range("A72").Value = 66
$ws.Range("A73").Value = 67
